$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the added columns
$ws.Range("Y1").Value = "cdc_revision"
$ws.Range("Z1").Value = "cdc_nice_name"
$ws.Range("AA1").Value = "reviewer_remove"

# Row 2
$ws.Range("Y2").Value = "not_applicable_germline"
$ws.Range("Z2").Value = "Not Applicable Germline"
$ws.Range("AA2").Value = $false

# Row 3
$ws.Range("Y3").Value = "undifferentiated_pleomorphic_sarcoma"
$ws.Range("Z3").Value = "Undifferentiated Pleomorphic Sarcoma"
$ws.Range("AA3").Value = $false

# Row 4
$ws.Range("Y4").Value = "undifferentiated_pleomorphic_sarcoma"
$ws.Range("Z4").Value = "Undifferentiated Pleomorphic Sarcoma"
$ws.Range("AA4").Value = $false

# Row 5
$ws.Range("Y5").Value = "not_applicable_germline"
$ws.Range("Z5").Value = "Not Applicable Germline"
$ws.Range("AA5").Value = $false

# Row 6
$ws.Range("Y6").Value = "myxoid_chondrosarcoma"
$ws.Range("Z6").Value = "Myxoid Chondrosarcoma"
$ws.Range("AA6").Value = $false

# Match column width customization seen for column Z (26th column) in the target
# (value chosen so the stored OOXML width lands as close as possible to 31.1640625,
# matching the bestFit width Excel computed for "Not Applicable Germline")
$ws.Columns.Item(26).ColumnWidth = 30.3307291666667

# Reflect the new selection location noted in the sheet view
$ws.Range("AA8").Select()
